$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text-like values directly
$ws.Range("A2").Value = "age_gr=30-39"
$ws.Range("B2").Value = "-0.07***"
$ws.Range("C2").Value = "-0.07***"
$ws.Range("D2").Value = "-0.07***"
$ws.Range("E2").Value = "-0.16***"
$ws.Range("F2").Value = "-0.16***"
$ws.Range("G2").Value = "-0.15***"
$ws.Range("B3").Value = "(0.01)"
$ws.Range("C3").Value = "(0.01)"
$ws.Range("D3").Value = "(0.01)"
$ws.Range("E3").Value = "(0.01)"
$ws.Range("F3").Value = "(0.01)"
$ws.Range("G3").Value = "(0.01)"
$ws.Range("A4").Value = "age_gr=40-48"
$ws.Range("B4").Value = "-0.10***"
$ws.Range("C4").Value = "-0.10***"
$ws.Range("D4").Value = "-0.09***"
$ws.Range("E4").Value = "-0.24***"
$ws.Range("F4").Value = "-0.25***"
$ws.Range("G4").Value = "-0.23***"
$ws.Range("B5").Value = "(0.01)"
$ws.Range("C5").Value = "(0.01)"
$ws.Range("D5").Value = "(0.01)"
$ws.Range("E5").Value = "(0.01)"
$ws.Range("F5").Value = "(0.01)"
$ws.Range("G5").Value = "(0.01)"
$ws.Range("A6").Value = "age_gr=49-57"
$ws.Range("B6").Value = "-0.13***"
$ws.Range("C6").Value = "-0.13***"
$ws.Range("D6").Value = "-0.12***"
$ws.Range("E6").Value = "-0.28***"
$ws.Range("F6").Value = "-0.30***"
$ws.Range("G6").Value = "-0.27***"
$ws.Range("B7").Value = "(0.01)"
$ws.Range("C7").Value = "(0.01)"
$ws.Range("D7").Value = "(0.01)"
$ws.Range("E7").Value = "(0.02)"
$ws.Range("F7").Value = "(0.02)"
$ws.Range("G7").Value = "(0.02)"
$ws.Range("A8").Value = "age_gr=>57"
$ws.Range("B8").Value = "-0.08***"
$ws.Range("C8").Value = "-0.09***"
$ws.Range("D8").Value = "-0.08***"
$ws.Range("E8").Value = "-0.22***"
$ws.Range("F8").Value = "-0.23***"
$ws.Range("G8").Value = "-0.21***"
$ws.Range("B9").Value = "(0.01)"
$ws.Range("C9").Value = "(0.01)"
$ws.Range("D9").Value = "(0.01)"
$ws.Range("E9").Value = "(0.02)"
$ws.Range("F9").Value = "(0.02)"
$ws.Range("G9").Value = "(0.02)"
$ws.Range("A10").Value = "educ_gr=low educ"
$ws.Range("D10").Value = "-0.01**"
$ws.Range("F10").Value = "0.07***"
$ws.Range("G10").Value = "0.03***"
$ws.Range("C11").Value = "(0.00)"
$ws.Range("D11").Value = "(0.01)"
$ws.Range("F11").Value = "(0.01)"
$ws.Range("G11").Value = "(0.01)"
$ws.Range("A12").Value = "HHinc_gr=low inc"
$ws.Range("D12").Value = "0.06***"
$ws.Range("G12").Value = "0.18***"
$ws.Range("D13").Value = "(0.00)"
$ws.Range("G13").Value = "(0.01)"
$ws.Range("A14").Value = "expvol"
$ws.Range("B14").Value = "0.41***"
$ws.Range("C14").Value = "0.41***"
$ws.Range("D14").Value = "0.45***"
$ws.Range("E14").Value = "0.93***"
$ws.Range("F14").Value = "0.90***"
$ws.Range("G14").Value = "1.01***"
$ws.Range("B15").Value = "(0.16)"
$ws.Range("C15").Value = "(0.16)"
$ws.Range("D15").Value = "(0.16)"
$ws.Range("E15").Value = "(0.27)"
$ws.Range("F15").Value = "(0.27)"
$ws.Range("G15").Value = "(0.27)"
$ws.Range("A16").Value = "N"
$ws.Range("A17").Value = "R2"

# Set numeric-looking values while forcing text storage (to match shared string type)
$numCells = @("C10","B16","C16","D16","E16","F16","G16","B17","C17","D17","E17","F17","G17")
foreach ($addr in $numCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("C10").Value = "0.00"
$ws.Range("B16").Value = "40529"
$ws.Range("C16").Value = "40529"
$ws.Range("D16").Value = "40529"
$ws.Range("E16").Value = "44874"
$ws.Range("F16").Value = "44874"
$ws.Range("G16").Value = "44874"
$ws.Range("B17").Value = "0.01"
$ws.Range("C17").Value = "0.01"
$ws.Range("D17").Value = "0.01"
$ws.Range("E17").Value = "0.01"
$ws.Range("F17").Value = "0.02"
$ws.Range("G17").Value = "0.03"
foreach ($addr in $numCells) {
    $ws.Range($addr).Style = "Normal"
}
